$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New incident rows appended to the "reportes" sheet (rows 4-7)
$data = @(
    @(3, "Tráfico pesado", "eeeeee",           4.873884452064224,  -74.03782881343682, "2025-11-13 17:12:27", 0),
    @(4, "Obstrucción",    "rtrtrtr",           4.873369135084946,  -74.03789817478737, "2025-11-13 17:14:35", 0),
    @(5, "Tráfico pesado", "ddsdasas",          4.87304,             -74.03791,          "2025-11-13 17:25:10", 1),
    @(6, "Accidente",      "1233333321332231",  4.873106115421405,  -74.03796509618449, "2025-11-13 17:27:10", 1)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Description column: values are always plain text in the source data.
    # A purely numeric-looking description (row 7) must be forced to text so
    # it isn't reinterpreted as a number by Excel.
    $cText = $row[2]
    $cCell = $ws.Cells.Item($r, 3)
    if ($cText -match '^-?\d+$') {
        $cCell.NumberFormat = "@"
    }
    $cCell.Value = $cText

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
